$wb = $excel.ActiveWorkbook
$sheet2 = $wb.Worksheets.Item("Sheet2")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$newSheet.Name = "Planilha2"

$src = $sheet2.Range("A12:C20")
$src.Copy($newSheet.Range("A1"))

$newSheet.Range("C9").Formula = "=SUM(C2:C8)"
$newSheet.Range("C3").Formula = "=B3-B2"
$newSheet.Range("C4").Formula = "=B4-B3"
$newSheet.Range("C5").Formula = "=B5-B4"
$newSheet.Range("C6").Formula = "=B6-B5"
$newSheet.Range("C7").Formula = "=B7-B6"
$newSheet.Range("C8").Formula = "=B8-B7"
